# Financials update: insert a new reporting-period column (column D) in front of
# the existing D:K data, shifting the old columns right by one (D->E, E->F, ... J->K),
# and refresh the historical figures with the restated numbers from the updated
# source feed. The brand-new column D carries the newest period (31-Dec-18,
# serial 43465) while the rightmost shifted column (K) keeps the oldest period's
# original values untouched (old J == new K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new column at D; Excel shifts D:K -> E:L, carrying values/number
#    formats with them (so the old "K7 empty / style-2" date cell now becomes
#    L7, and so on for every row).
$ws.Columns.Item(4).Insert()

# 2) The freshly inserted column D has no number format of its own yet (Excel
#    defaults it to the sheet's base style). Clone the format from column E
#    (which now holds what used to be column D) back onto D, one contiguous
#    block of rows at a time so we don't accidentally stamp a style onto rows
#    that never had D:K data (the plain label rows 5, 6, 37, 79, ...).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Write the refreshed figures. Keys are 1-based column numbers
#    (4=D, 5=E, 6=F, 7=G, 8=H, 9=I, 10=J, 11=K); column L is left as the
#    (already shifted/formatted) blank cell Excel produced during the insert.
$data = @{
    7 = @{4=43465; 5=43100; 6=42735; 7=42369; 8=42004; 9=41639; 10=41274; 11=40908}
    8 = @{4=1450100; 5=1345300; 6=1192700; 7=1123800; 8=1041500; 9=979200; 10=848300; 11=855000}
    9 = @{4=295900; 5=254400; 6=208300; 7=192300; 8=184700; 9=175200; 10=158300; 11=158700}
    10 = @{4=1154200; 5=1090900; 6=984400; 7=931500; 8=856800; 9=804000; 10=690000; 11=696300}
    12 = @{4=67300; 5=70400; 6=47600; 7=34100; 8=26400; 9=37800; 10=34200; 11=30800}
    13 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    14 = @{4=3000; 5=400; 6=3400; 7=6800; 8=2100; 9=7700; 10=-31100; 11=-42100}
    15 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    17 = @{4=1129200; 5=1016200; 6=885500; 7=850600; 8=739300; 9=720600; 10=614600; 11=579000}
    18 = @{4=320900; 5=329000; 6=307200; 7=273300; 8=302200; 9=258500; 10=233700; 11=276000}
    20 = @{4=15300; 5=6300; 6=15700; 7=34500; 8=23700; 9=30900; 10=-13000; 11=-22100}
    21 = @{4=480200; 5=461600; 6=424200; 7=382400; 8=394700; 9=354900; 10=288500; 11=319000}
    22 = @{4=0; 5=0; 6=100; 7=0; 8=200; 9=100; 10=0; 11=0}
    23 = @{4=336200; 5=335300; 6=322900; 7=307700; 8=325700; 9=289300; 10=220600; 11=253900}
    24 = @{4=80400; 5=103100; 6=100000; 7=113900; 8=124000; 9=112200; 10=99100; 11=100000}
    25 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    26 = @{4=255800; 5=232300; 6=222900; 7=193800; 8=201600; 9=177100; 10=121600; 11=153900}
    27 = @{4=256000; 5=232200; 6=222800; 7=193800; 8=201600; 9=177100; 10=121600; 11=153800}
    28 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    29 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    30 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    31 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    32 = @{4=-15300; 5=-6300; 6=-15700; 7=-34500; 8=-23700; 9=-30900; 10=13000; 11=22100}
    33 = @{4=256000; 5=232200; 6=222800; 7=193800; 8=201600; 9=177100; 10=121600; 11=153800}
    34 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    35 = @{4=256000; 5=232200; 6=222800; 7=193800; 8=201600; 9=177100; 10=121600; 11=153800}
    38 = @{4=43465; 5=43100; 6=42735; 7=42369; 8=42004; 9=41639; 10=41274; 11=40908}
    41 = @{4=1053100; 5=776500; 6=764000; 7=667500; 8=570500; 9=719900; 10=494800; 11=637600}
    42 = @{4=554100; 5=601500; 6=453400; 7=650100; 8=671900; 9=551100; 10=579600; 11=334700}
    43 = @{4=352400; 5=359600; 6=321400; 7=270900; 8=235700; 9=219900; 10=188900; 11=185900}
    44 = @{4=33200; 5=26900; 6=15500; 7=7300; 8=4700; 9=3700; 10=3300; 11=3800}
    45 = @{4=194100; 5=195300; 6=208900; 7=225600; 8=178700; 9=174300; 10=173100; 11=179300}
    46 = @{4=2186900; 5=1959800; 6=1763100; 7=1821500; 8=1661500; 9=1668900; 10=1439700; 11=1341300}
    47 = @{4=404400; 5=473600; 6=458600; 7=567700; 8=624200; 9=462200; 10=320800; 11=208300}
    48 = @{4=79300; 5=74300; 6=59200; 7=52300; 8=44700; 9=41400; 10=40100; 11=45700}
    49 = @{4=317400; 5=364300; 6=384500; 7=74200; 8=85200; 9=82800; 10=71100; 11=84400}
    50 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    51 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    52 = @{4=142000; 5=121600; 6=123800; 7=110600; 8=115000; 9=108500; 10=108100; 11=109900}
    53 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    54 = @{4=3130000; 5=2993700; 6=2789200; 7=2626300; 8=2530600; 9=2363900; 10=1979800; 11=1789600}
    57 = @{4=14700; 5=6600; 6=8100; 7=5100; 8=6100; 9=4400; 10=6600; 11=7900}
    58 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    59 = @{4=994100; 5=976600; 6=911200; 7=849600; 8=832700; 9=782700; 10=698900; 11=625500}
    60 = @{4=1008800; 5=983200; 6=919300; 7=854600; 8=838700; 9=787100; 10=705500; 11=633300}
    61 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    62 = @{4=430000; 5=409700; 6=365000; 7=328000; 8=307900; 9=288300; 10=248300; 11=204000}
    63 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    64 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    65 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    66 = @{4=1439500; 5=1392900; 6=1284300; 7=1182700; 8=1146800; 9=1075400; 10=953900; 11=837800}
    68 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    69 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    70 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    71 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    72 = @{4=1352600; 5=1287800; 6=1229900; 7=1143900; 8=1098500; 9=1047100; 10=964200; 11=889800}
    73 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    74 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    75 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    76 = @{4=1690500; 5=1600700; 6=1504900; 7=1443600; 8=1383900; 9=1288500; 10=1025900; 11=951800}
    77 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    80 = @{4=43465; 5=43100; 6=42735; 7=42369; 8=42004; 9=41639; 10=41274; 11=40908}
    81 = @{4=256000; 5=232200; 6=222800; 7=193800; 8=201600; 9=177100; 10=121600; 11=153800}
    83 = @{4=143900; 5=126100; 6=101200; 7=74600; 8=68800; 9=65500; 10=67800; 11=65100}
    84 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    85 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    86 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    87 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    88 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    89 = @{4=451600; 5=424100; 6=302900; 7=275600; 8=288800; 9=277900; 10=285000; 11=231700}
    91 = @{4=-32200; 5=-41600; 6=-27000; 7=-31900; 8=-18400; 9=-16600; 10=-10800; 11=-61900}
    92 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    93 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    94 = @{4=7400; 5=-305700; 6=116800; 7=-44500; 8=-297600; 9=-146000; 10=-412600; 11=-23200}
    96 = @{4=-178900; 5=-168900; 6=-131000; 7=-134100; 8=-142600; 9=-76100; 10=-97600; 11=-79200}
    97 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    98 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    99 = @{4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0}
    100 = @{4=-119500; 5=-152800; 6=-136100; 7=-84300; 8=-152700; 9=-1900; 10=-100600; 11=-120300}
    101 = @{4=-39500; 5=-1100; 6=2800; 7=-38200; 8=23300; 9=82500; 10=40900; 11=-25800}
    102 = @{4=300000; 5=-35600; 6=286500; 7=108600; 8=-138200; 9=212500; 10=-187300; 11=62400}
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($c in $rowVals.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
